# Update "想去人数" (want-to-go count) values in column F across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1124
$ws1.Range("F5").Value = 182
$ws1.Range("F7").Value = 230
$ws1.Range("F11").Value = 515
$ws1.Range("F14").Value = 12791
$ws1.Range("F15").Value = 7
$ws1.Range("F16").Value = 5257
$ws1.Range("F17").Value = 5530

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 105

# Sheet "全部类型" (All types, aggregate)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1124
$ws4.Range("F6").Value = 182
$ws4.Range("F8").Value = 230
$ws4.Range("F12").Value = 515
$ws4.Range("F15").Value = 12791
$ws4.Range("F16").Value = 105
$ws4.Range("F18").Value = 7
$ws4.Range("F19").Value = 5257
$ws4.Range("F20").Value = 5530
